$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2218649517684887
$ws.Range("C2").Value = 0.4919614147909968
$ws.Range("J2").Value = 0.01929260450160772
$ws.Range("P2").Value = 0.1607717041800643
$ws.Range("S2").Value = 0.1061093247588424

# Row 3
$ws.Range("B3").Value = 0.0245398773006135
$ws.Range("C3").Value = 0.05521472392638037
$ws.Range("J3").Value = 0.049079754601227
$ws.Range("P3").Value = 0.7116564417177914
$ws.Range("S3").Value = 0.1595092024539877

# Row 4
$ws.Range("J4").Value = 0.05405405405405406
$ws.Range("P4").Value = 0.7027027027027027
$ws.Range("S4").Value = 0.2432432432432433

# Row 6
$ws.Range("B6").Value = 0.04186046511627907
$ws.Range("D6").Value = 0.009302325581395349
$ws.Range("F6").Value = 0.04651162790697674
$ws.Range("J6").Value = 0.2186046511627907
$ws.Range("O6").Value = 0.009302325581395349
$ws.Range("Q6").Value = 0.1813953488372093
$ws.Range("R6").Value = 0.06511627906976744
$ws.Range("S6").Value = 0.427906976744186

# Row 7
$ws.Range("B7").Value = 0.1074766355140187
$ws.Range("D7").Value = 0.01401869158878505
$ws.Range("F7").Value = 0.04205607476635514
$ws.Range("J7").Value = 0.1869158878504673
$ws.Range("O7").Value = 0.009345794392523364
$ws.Range("Q7").Value = 0.1775700934579439
$ws.Range("R7").Value = 0.05607476635514019
$ws.Range("S7").Value = 0.4065420560747663

# Row 8
$ws.Range("B8").Value = 0.0752212389380531
$ws.Range("D8").Value = 0.00663716814159292
$ws.Range("F8").Value = 0.06858407079646017
$ws.Range("J8").Value = 0.1349557522123894
$ws.Range("O8").Value = 0.008849557522123894
$ws.Range("Q8").Value = 0.163716814159292
$ws.Range("R8").Value = 0.08628318584070796
$ws.Range("S8").Value = 0.4557522123893805

# Row 9
$ws.Range("B9").Value = 0.05
$ws.Range("D9").Value = 0.025
$ws.Range("F9").Value = 0.035
$ws.Range("J9").Value = 0.1
$ws.Range("O9").Value = 0.015
$ws.Range("Q9").Value = 0.16
$ws.Range("R9").Value = 0.04
$ws.Range("S9").Value = 0.575

# Row 10
$ws.Range("B10").Value = 0.1092553931802366
$ws.Range("D10").Value = 0.0173973556019485
$ws.Range("E10").Value = 0.002783576896311761
$ws.Range("F10").Value = 0.07306889352818371
$ws.Range("J10").Value = 0.1356993736951983
$ws.Range("O10").Value = 0.01322199025748086
$ws.Range("Q10").Value = 0.2066805845511482
$ws.Range("R10").Value = 0.05984690327070286
$ws.Range("S10").Value = 0.3820459290187891

# Row 11
$ws.Range("F11").Value = 0.003003003003003003
$ws.Range("G11").Value = 0.1561561561561562
$ws.Range("J11").Value = 0.09309309309309309
$ws.Range("K11").Value = 0.1861861861861862
$ws.Range("L11").Value = 0.5555555555555556
$ws.Range("S11").Value = 0.006006006006006006

# Row 12
$ws.Range("G12").Value = 0.7150259067357513
$ws.Range("J12").Value = 0.2227979274611399
$ws.Range("L12").Value = 0.02072538860103627
$ws.Range("S12").Value = 0.04145077720207254

# Row 13
$ws.Range("G13").Value = 0.7105263157894737
$ws.Range("J13").Value = 0.2894736842105263

# Row 15
$ws.Range("F15").Value = 0.01754385964912281
$ws.Range("H15").Value = 0.1535087719298246
$ws.Range("I15").Value = 0.07456140350877193
$ws.Range("J15").Value = 0.3947368421052632
$ws.Range("K15").Value = 0.07894736842105263
$ws.Range("M15").Value = 0.008771929824561403
$ws.Range("O15").Value = 0.04385964912280702
$ws.Range("S15").Value = 0.2280701754385965

# Row 16
$ws.Range("F16").Value = 0.01081081081081081
$ws.Range("H16").Value = 0.1351351351351351
$ws.Range("I16").Value = 0.0918918918918919
$ws.Range("J16").Value = 0.4540540540540541
$ws.Range("K16").Value = 0.1189189189189189
$ws.Range("M16").Value = 0.03243243243243243
$ws.Range("O16").Value = 0.02702702702702703
$ws.Range("S16").Value = 0.1297297297297297

# Row 17
$ws.Range("F17").Value = 0.01902748414376321
$ws.Range("H17").Value = 0.1818181818181818
$ws.Range("I17").Value = 0.08456659619450317
$ws.Range("J17").Value = 0.4820295983086681
$ws.Range("K17").Value = 0.07822410147991543
$ws.Range("M17").Value = 0.008456659619450317
$ws.Range("O17").Value = 0.05496828752642706
$ws.Range("S17").Value = 0.09090909090909091

# Row 18
$ws.Range("F18").Value = 0.01257861635220126
$ws.Range("H18").Value = 0.1761006289308176
$ws.Range("I18").Value = 0.1069182389937107
$ws.Range("J18").Value = 0.389937106918239
$ws.Range("K18").Value = 0.1069182389937107
$ws.Range("M18").Value = 0.01886792452830189
$ws.Range("N18").Value = 0.006289308176100629
$ws.Range("O18").Value = 0.06289308176100629
$ws.Range("S18").Value = 0.119496855345912

# Row 19
$ws.Range("F19").Value = 0.01103752759381899
$ws.Range("H19").Value = 0.2060338484179544
$ws.Range("I19").Value = 0.08094186902133922
$ws.Range("J19").Value = 0.3899926416482708
$ws.Range("K19").Value = 0.1250919793966151
$ws.Range("M19").Value = 0.01766004415011038
$ws.Range("O19").Value = 0.08462104488594555
$ws.Range("S19").Value = 0.08462104488594555
